# Generate Report for Handback
#
# Marks the zh-cn and de-de localization rows as handed back: updates the
# shared "Status" text (used on the Overview sheet and in each language
# table), fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for each language, adds a hyperlink
# from the new "Latest Target File" cell to the source markdown file, and
# widens a few columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdDisplay = "84191805-2b27-4f11-8519-64ea8ae9ae5f.md"
$mdTarget  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbfc818b29d5655fa75dcdef76cbf44c62203415/e2e/84191805-2b27-4f11-8519-64ea8ae9ae5f.md"

# --- Overview sheet: the "zh-cn" / "de-de" status columns share the same
#     underlying text as each language table's "Status" column. ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.1667
$wsZh.Columns.Item(10).ColumnWidth = 39.1667

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdTarget, "", "", $mdDisplay)
$wsZh.Range("J2").Value = "84191805-2b27-4f11-8519-64ea8ae9ae5f.5820819b002e23390a75775157cee49df6858b01.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-26 20:57:28"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.1667
$wsDe.Columns.Item(10).ColumnWidth = 39.1667

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdTarget, "", "", $mdDisplay)
$wsDe.Range("J2").Value = "84191805-2b27-4f11-8519-64ea8ae9ae5f.5820819b002e23390a75775157cee49df6858b01.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-26 20:57:35"
